$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.704.17'
$ws.Range("E2").Value = '  +5.71%  '

$ws.Range("D3").Value = '3.241.15'
$ws.Range("E3").Value = '  +2.51%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.32%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -1.78%  '

$ws.Range("D9").Value = '3.236.78'
$ws.Range("E9").Value = '  +2.46%  '

$ws.Range("E10").Value = '  +4.87%  '

$ws.Range("E11").Value = '  +3.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.413'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.07%  '

$ws.Range("D13").Value = '3.804.17'
$ws.Range("E13").Value = '  +2.50%  '

$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.16%  '

$ws.Range("D16").Value = '67.624.97'
$ws.Range("E16").Value = '  +5.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000168'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.22%  '

$ws.Range("D18").Value = '3.238.70'
$ws.Range("E18").Value = '  +2.57%  '

$ws.Range("E19").Value = '  +1.57%  '

$ws.Range("E20").Value = '  +3.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.48%  '

$ws.Range("E22").Value = '  +5.58%  '

$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.68%  '

$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("E26").Value = '  +3.87%  '

$ws.Range("E27").Value = '  +1.78%  '

$ws.Range("E28").Value = '  +3.32%  '

$ws.Range("E29").Value = '  +0.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.89%  '

$ws.Range("E31").Value = '  +3.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.30%  '

$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.90%  '

$ws.Range("E37").Value = '  +5.07%  '

$ws.Range("E38").Value = '  +5.73%  '

$ws.Range("E39").Value = '  +10.65%  '

$ws.Range("E40").Value = '  +15.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.07%  '

$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '363.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.67%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.17%  '

$ws.Range("E44").Value = '  +6.23%  '

$ws.Range("D45").Value = '2.698.05'
$ws.Range("E45").Value = '  +1.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.93%  '

$ws.Range("E48").Value = '  +3.85%  '

$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.03%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.103'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.35%  '

